# Update cryptos list: price (D) and volume% (E) columns for rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.750.23'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.630.81'
$ws.Range("E3").Value = '  -0.14%  '
$ws.Range("E4").Value = '  -0.62%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.500'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  -0.99%  '
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0792'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("E12").Value = '  +0.40%  '
$ws.Range("D13").Value = '1.855.52'
$ws.Range("E13").Value = '  -0.22%  '
$ws.Range("D14").Value = '1.628.63'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("E15").Value = '  +0.28%  '
$ws.Range("D16").Value = '0.0₃0760'
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("D18").Value = '25.745.10'
$ws.Range("E18").Value = '  -0.22%  '
$ws.Range("E19").Value = '  -0.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '191.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("E22").Value = '  -0.23%  '
$ws.Range("E23").Value = '  +1.17%  '
$ws.Range("E24").Value = '  +2.35%  '
$ws.Range("E25").Value = '  -0.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '142.23'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.79%  '
$ws.Range("E27").Value = '  +2.81%  '
$ws.Range("E28").Value = '  +0.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("E30").Value = '  -0.75%  '
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("E33").Value = '  -0.70%  '
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.905'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.00%  '
$ws.Range("D37").Value = '1.138.42'
$ws.Range("E37").Value = '  +2.81%  '
$ws.Range("E38").Value = '  -2.76%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.542'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("E40").Value = '  -0.55%  '
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("E42").Value = '  -0.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.13%  '
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.806'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.57%  '
$ws.Range("D46").Value = '1.765.24'
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.19'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0510'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.44'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.36'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.13%  '
